$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = 0
$ws.Range("C27").Value = 0

$ws.Range("E13").Select()
